# Update capital structure database values for rows 2 and 3
# (Romania / Brokerage & Investment Banking rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> new value, applied identically to rows 2 and 3
$updates = [ordered]@{
    "D"  = -0.27
    "K"  = -1.82
    "L"  = -1.913774973711882
    "M"  = 0.258
    "N"  = 0.03712230215827338
    "O"  = -0.1417582417582418
    "P"  = 0.258
    "Q"  = 0.03712230215827338
    "R"  = -0.1417582417582418
    "T"  = 0
    "U"  = 1.9
    "V"  = 0.2733812949640287
    "W"  = -0.1166666666666667
    "X"  = 0.03902470745751097
    "Y"  = -0.1556913741241776
    "Z"  = 0.07342495367510808
    "AB" = 0.03872579729654586
    "AC" = -0.03872579729654586
    "AD" = 0.332
    "AF" = 0.332
    "AG" = -1.568
    "AH" = 0.04559187036528427
    "AI" = 0.02400231347599769
    "AJ" = -0.291341508732813
    "AK" = -0.1314113308749581
}

foreach ($row in 2,3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
